$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Test Two"
$ws.Range("B4").Value = "Testing form submission to Excel file                            "
$ws.Range("C4").Value = "test@two.com"
$ws.Range("D4").Value = "Urgent"
